# Generate Report for Handoff
# Refreshes the "Latest Handoff Datetime" values that were generated for this
# handoff run. Rows whose datetime matched the previous batch timestamp get
# updated to the new timestamp produced by this report generation; rows that
# already have their own distinct datetime (e.g. files still "In Translation")
# are left untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value  = "2016-03-23 10:27:00"
$wsOverview.Range("D6").Value  = "2016-03-23 10:27:00"
$wsOverview.Range("D7").Value  = "2016-03-23 10:27:00"
$wsOverview.Range("D8").Value  = "2016-03-23 10:27:00"
$wsOverview.Range("D9").Value  = "2016-03-23 10:27:00"
$wsOverview.Range("D10").Value = "2016-03-23 10:27:00"

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value  = "2016-03-23 10:26:52"
$wsZhCn.Range("E6").Value  = "2016-03-23 10:26:52"
$wsZhCn.Range("E7").Value  = "2016-03-23 10:26:52"
$wsZhCn.Range("E8").Value  = "2016-03-23 10:26:52"
$wsZhCn.Range("E9").Value  = "2016-03-23 10:26:52"
$wsZhCn.Range("E10").Value = "2016-03-23 10:26:52"

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value  = "2016-03-23 10:27:00"
$wsDeDe.Range("E6").Value  = "2016-03-23 10:27:00"
$wsDeDe.Range("E7").Value  = "2016-03-23 10:27:00"
$wsDeDe.Range("E8").Value  = "2016-03-23 10:27:00"
$wsDeDe.Range("E9").Value  = "2016-03-23 10:27:00"
$wsDeDe.Range("E10").Value = "2016-03-23 10:27:00"
